# Apply the LinuxForHealth rebrand + version bump edit to the
# StructureDefinition-insight-detail workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "Metadata" sheet - top of file summary properties
# ---------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-detail"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# ---------------------------------------------------------------
# "Elements" sheet - structure definition element table
# ---------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root "Extension" row no longer repeats the ele-1/ext-1 constraint text
$elements.Range("AI2").Value = ""

# Slice rows: extension type references point at the new host
$elements.Range("J5").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference}" + [char]10
$elements.Range("J6").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/reference-path}" + [char]10
$elements.Range("J7").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/evaluated-output}" + [char]10
$elements.Range("J8").Value = "Extension {http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-result}" + [char]10

# Extension.url fixed value
$elements.Range("Q9").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/insight-detail"
